# daily auto push: append the next day's row (2025/10/01, 水, 0, 142)
# to the bottom of the data table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 42

# Column A holds the date as literal text (e.g. "2025/09/30"), not a real
# Excel date serial. Writing a date-shaped string via .Value would normally
# be auto-recognized as a date and stamped with a date NumberFormat, so we
# format the cell as Text first and clear the resulting format afterwards -
# this keeps the stored value a plain string with no explicit cell style,
# matching every other row in the table.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025/10/01"
$ws.Range("A" + $newRow).ClearFormats()

$ws.Range("B" + $newRow).Value = "水"
$ws.Range("C" + $newRow).Value = 0
$ws.Range("D" + $newRow).Value = 142
